$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G ("K") values per row, replacing the old Strike# derived values
# with the newly regenerated K counts.
$gValues = @{
    2 = 1
    4 = 1
    5 = 0
    6 = 1
    7 = 1
    8 = 0
    9 = 1
    10 = 0
    11 = 0
    12 = 1
    13 = 2
    14 = 1
    15 = 1
    16 = 1
    17 = 1
    18 = 1
    19 = 0
    20 = 0
    21 = 1
    22 = 1
    23 = 0
    24 = 1
    25 = 0
    26 = 2
    27 = 2
    28 = 0
    29 = 1
    30 = 2
    31 = 2
    32 = 3
    33 = 2
    34 = 1
    35 = 2
    36 = 2
    38 = 2
    39 = 1
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}
